# Apply an AutoFilter over the full data range (A1:D50) on Sheet1, which is
# the core change behind this edit (drives the <autoFilter> element on the
# sheet plus the hidden _xlnm._FilterDatabase defined name scoped to the
# sheet in the workbook part).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filterRange = $ws.Range("A1:D50")
$filterRange.AutoFilter()

# Excel records the filtered range as a workbook-level, sheet-scoped, hidden
# defined name called _xlnm._FilterDatabase.
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$50")
$filterDbName.Visible = $false

# Column A ("Vehicle Type" text) and column B ("Model Year") were resized to
# fit their contents once the filter drop-downs were added.
$ws.Columns.Item(1).ColumnWidth = 14.917
$ws.Columns.Item(2).ColumnWidth = 9.25

# The active cell ended up on H4 after the interaction that produced this
# workbook state.
$ws.Range("H4").Select()
